# Update view-count figures ("F" column) on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets to reflect refreshed scrape data.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F15").Value = 23
$wsExpo.Range("F16").Value = 6522
$wsExpo.Range("F19").Value = 138
$wsExpo.Range("F22").Value = 15582
$wsExpo.Range("F26").Value = 106

# Sheet 4: 全部类型
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F18").Value = 23
$wsAll.Range("F19").Value = 6522
$wsAll.Range("F22").Value = 138
$wsAll.Range("F26").Value = 15582
$wsAll.Range("F30").Value = 106
